$d = $word.ActiveDocument

function Find-ParaIndex($pattern) {
    $idx = 0
    $i = 0
    foreach ($p in $d.Paragraphs) {
        $i = $i + 1
        if ($idx -eq 0 -and $p.Range.Text -like $pattern) {
            $idx = $i
        }
    }
    return $idx
}

# 1) The "Bases de données" paragraph becomes "Visualisation : tableau",
#    and a brand new "MLOps : ..." paragraph is inserted right after it.
$basesIndex = Find-ParaIndex("Bases de donn*es : SQL, MongoDB, Neo4j, Redis*")
$basesPara = $d.Paragraphs($basesIndex)
$basesPara.Range.Text = "Visualisation : tableau"
$basesPara.Range.InsertParagraphAfter()
$d.Paragraphs($basesIndex + 1).Range.Text = "MLOps : spark, Git, DVC, Flask, Docker, Github Actions, Heroku, MLflow, Streamlit"

# 2) The old "Visualisation : tableau" paragraph (the second one with that
#    text now present in the document) is removed entirely.
$visuIndex = 0
$i = 0
foreach ($p in $d.Paragraphs) {
    $i = $i + 1
    if ($p.Range.Text -like "Visualisation : tableau*") {
        $visuIndex = $i
    }
}
$d.Paragraphs($visuIndex).Range.Delete()

# 3) The old "MLOps : ..." paragraph (the second one with that text now
#    present in the document) becomes "Bases de données : SQL, MongoDB,
#    Neo4j, Redis".
$mlopsIndex = 0
$i = 0
foreach ($p in $d.Paragraphs) {
    $i = $i + 1
    if ($p.Range.Text -like "MLOps : spark, Git, DVC, Flask, Docker, Github Actions, Heroku, MLflow, Streamlit*") {
        $mlopsIndex = $i
    }
}
$d.Paragraphs($mlopsIndex).Range.Text = "Bases de données : SQL, MongoDB, Neo4j, Redis"
